$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on humidity-percentage cells first so Excel does not
# auto-convert the literal "NN%" strings into numeric percentage values.
$pctCells = @("H3", "H6", "H7", "H8", "H13", "H18", "H22", "H27", "H29", "H32", "H33", "H39", "H40", "H46")
foreach ($c in $pctCells) {
    $ws.Range($c).NumberFormat = "@"
}

# Apply the updated extraction timestamps and measurement readings.
$ws.Range("E2").Value = "2026-02-27 17:18:16"
$ws.Range("K2").Value = "11.2 MJ/m2"
$ws.Range("O2").Value = "5.8 °C"
$ws.Range("E3").Value = "2026-02-27 17:18:18"
$ws.Range("H3").Value = "36%"
$ws.Range("K3").Value = "18.8 MJ/m2"
$ws.Range("E4").Value = "2026-02-27 17:18:21"
$ws.Range("J4").Value = "1024.8 hPa"
$ws.Range("K4").Value = "8.8 MJ/m2"
$ws.Range("O4").Value = "9.8 °C"
$ws.Range("E5").Value = "2026-02-27 17:18:23"
$ws.Range("K5").Value = "16.1 MJ/m2"
$ws.Range("O5").Value = "5.3 °C"
$ws.Range("E6").Value = "2026-02-27 17:18:26"
$ws.Range("H6").Value = "88%"
$ws.Range("J6").Value = "1024.7 hPa"
$ws.Range("K6").Value = "9.4 MJ/m2"
$ws.Range("O6").Value = "11.1 °C"
$ws.Range("E7").Value = "2026-02-27 17:18:28"
$ws.Range("H7").Value = "87%"
$ws.Range("J7").Value = "1025.1 hPa"
$ws.Range("K7").Value = "13.8 MJ/m2"
$ws.Range("O7").Value = "11.4 °C"
$ws.Range("E8").Value = "2026-02-27 17:18:31"
$ws.Range("H8").Value = "58%"
$ws.Range("J8").Value = "1024.4 hPa"
$ws.Range("K8").Value = "14.9 MJ/m2"
$ws.Range("N8").Value = "10.2 °C 16:38 TU"
$ws.Range("E9").Value = "2026-02-27 17:18:33"
$ws.Range("K9").Value = "7.9 MJ/m2"
$ws.Range("E10").Value = "2026-02-27 17:18:35"
$ws.Range("K10").Value = "14.2 MJ/m2"
$ws.Range("E11").Value = "2026-02-27 17:18:38"
$ws.Range("O11").Value = "8.6 °C"
$ws.Range("E12").Value = "2026-02-27 17:18:40"
$ws.Range("E13").Value = "2026-02-27 17:18:42"
$ws.Range("H13").Value = "67%"
$ws.Range("J13").Value = "1026.2 hPa"
$ws.Range("K13").Value = "15.0 MJ/m2"
$ws.Range("O13").Value = "6.2 °C"
$ws.Range("E14").Value = "2026-02-27 17:18:45"
$ws.Range("K14").Value = "11.5 MJ/m2"
$ws.Range("O14").Value = "10.8 °C"
$ws.Range("E15").Value = "2026-02-27 17:18:47"
$ws.Range("O15").Value = "10.8 °C"
$ws.Range("E16").Value = "2026-02-27 17:18:49"
$ws.Range("K16").Value = "15.7 MJ/m2"
$ws.Range("N16").Value = "0.2 °C 16:58 TU"
$ws.Range("E17").Value = "2026-02-27 17:18:52"
$ws.Range("K17").Value = "16.5 MJ/m2"
$ws.Range("N17").Value = "6.0 °C 16:46 TU"
$ws.Range("O17").Value = "8.0 °C"
$ws.Range("E18").Value = "2026-02-27 17:18:54"
$ws.Range("H18").Value = "80%"
$ws.Range("E19").Value = "2026-02-27 17:18:57"
$ws.Range("K19").Value = "14.5 MJ/m2"
$ws.Range("E20").Value = "2026-02-27 17:18:59"
$ws.Range("K20").Value = "16.8 MJ/m2"
$ws.Range("E21").Value = "2026-02-27 17:19:01"
$ws.Range("J21").Value = "1024.9 hPa"
$ws.Range("K21").Value = "15.4 MJ/m2"
$ws.Range("O21").Value = "9.3 °C"
$ws.Range("E22").Value = "2026-02-27 17:19:04"
$ws.Range("H22").Value = "49%"
$ws.Range("K22").Value = "17.4 MJ/m2"
$ws.Range("E23").Value = "2026-02-27 17:19:06"
$ws.Range("K23").Value = "18.1 MJ/m2"
$ws.Range("E24").Value = "2026-02-27 17:19:09"
$ws.Range("J24").Value = "1024.0 hPa"
$ws.Range("K24").Value = "15.6 MJ/m2"
$ws.Range("O24").Value = "10.1 °C"
$ws.Range("E25").Value = "2026-02-27 17:19:11"
$ws.Range("K25").Value = "17.1 MJ/m2"
$ws.Range("E26").Value = "2026-02-27 17:19:13"
$ws.Range("K26").Value = "16.0 MJ/m2"
$ws.Range("E27").Value = "2026-02-27 17:19:16"
$ws.Range("H27").Value = "40%"
$ws.Range("K27").Value = "16.3 MJ/m2"
$ws.Range("E28").Value = "2026-02-27 17:19:18"
$ws.Range("J28").Value = "1025.0 hPa"
$ws.Range("K28").Value = "7.0 MJ/m2"
$ws.Range("O28").Value = "7.7 °C"
$ws.Range("E29").Value = "2026-02-27 17:19:21"
$ws.Range("H29").Value = "87%"
$ws.Range("K29").Value = "10.8 MJ/m2"
$ws.Range("E30").Value = "2026-02-27 17:19:23"
$ws.Range("J30").Value = "1024.8 hPa"
$ws.Range("K30").Value = "5.3 MJ/m2"
$ws.Range("E31").Value = "2026-02-27 17:19:26"
$ws.Range("J31").Value = "1024.4 hPa"
$ws.Range("E32").Value = "2026-02-27 17:19:28"
$ws.Range("H32").Value = "62%"
$ws.Range("K32").Value = "14.2 MJ/m2"
$ws.Range("O32").Value = "8.0 °C"
$ws.Range("E33").Value = "2026-02-27 17:19:31"
$ws.Range("H33").Value = "53%"
$ws.Range("J33").Value = "1024.4 hPa"
$ws.Range("K33").Value = "15.7 MJ/m2"
$ws.Range("O33").Value = "8.3 °C"
$ws.Range("E34").Value = "2026-02-27 17:19:33"
$ws.Range("K34").Value = "14.2 MJ/m2"
$ws.Range("O34").Value = "5.0 °C"
$ws.Range("E35").Value = "2026-02-27 17:19:36"
$ws.Range("J35").Value = "1022.8 hPa"
$ws.Range("K35").Value = "15.3 MJ/m2"
$ws.Range("O35").Value = "12.4 °C"
$ws.Range("E36").Value = "2026-02-27 17:19:38"
$ws.Range("J36").Value = "1025.1 hPa"
$ws.Range("K36").Value = "9.2 MJ/m2"
$ws.Range("E37").Value = "2026-02-27 17:19:41"
$ws.Range("J37").Value = "1025.3 hPa"
$ws.Range("O37").Value = "8.0 °C"
$ws.Range("E38").Value = "2026-02-27 17:19:43"
$ws.Range("K38").Value = "12.6 MJ/m2"
$ws.Range("O38").Value = "10.3 °C"
$ws.Range("E39").Value = "2026-02-27 17:19:45"
$ws.Range("H39").Value = "31%"
$ws.Range("K39").Value = "16.9 MJ/m2"
$ws.Range("E40").Value = "2026-02-27 17:19:48"
$ws.Range("H40").Value = "68%"
$ws.Range("J40").Value = "1025.4 hPa"
$ws.Range("O40").Value = "9.0 °C"
$ws.Range("E41").Value = "2026-02-27 17:19:50"
$ws.Range("J41").Value = "1024.9 hPa"
$ws.Range("K41").Value = "13.2 MJ/m2"
$ws.Range("E42").Value = "2026-02-27 17:19:53"
$ws.Range("E43").Value = "2026-02-27 17:19:55"
$ws.Range("K43").Value = "14.1 MJ/m2"
$ws.Range("O43").Value = "9.2 °C"
$ws.Range("E44").Value = "2026-02-27 17:19:57"
$ws.Range("K44").Value = "17.1 MJ/m2"
$ws.Range("E45").Value = "2026-02-27 17:20:00"
$ws.Range("J45").Value = "1022.0 hPa"
$ws.Range("E46").Value = "2026-02-27 17:20:02"
$ws.Range("H46").Value = "83%"
$ws.Range("J46").Value = "1024.5 hPa"
$ws.Range("K46").Value = "13.1 MJ/m2"
$ws.Range("L46").Value = "33.5 km/h - 174º 16:48 TU"
$ws.Range("O46").Value = "10.7 °C"
